$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Price column values are forced to text (matching the sheet's existing
# inline-string cell type) so numeric-looking prices like "0.999" are not
# reinterpreted as numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.280.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.510.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.78%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.507.55"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.136"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.99"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.096.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000206"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.38"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.501.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.40%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.277.95"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "426.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.597"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "78.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.632.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.28"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.36%  "
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.47"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.98%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.495.32"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.73"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.63"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.79"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "170.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.63%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0859"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.885"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "45.33"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("E47").Value = "  -6.93%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.51%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.74%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.76%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.946"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.81%  "
